$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-covered rows 2-10 (A:T) with new NATMI values

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.905615
$ws.Range("H2").Value = 35.716845
$ws.Range("I2").Value = 0.8197078149061106
$ws.Range("J2").Value = 0.8197078149061106
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1532026666666667
$ws.Range("N2").Value = 0.459608
$ws.Range("O2").Value = 0.03935316038222934
$ws.Range("P2").Value = 0.03935316038222934
$ws.Range("Q2").Value = 1.823971966306667
$ws.Range("R2").Value = 16.41574769676
$ws.Range("S2").Value = 0.03225809310656693
$ws.Range("T2").Value = 0.03225809310656693

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.905615
$ws.Range("H3").Value = 35.716845
$ws.Range("I3").Value = 0.8197078149061106
$ws.Range("J3").Value = 0.8197078149061106
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1693566666666667
$ws.Range("N3").Value = 0.50807
$ws.Range("O3").Value = 0.04350263745496
$ws.Range("P3").Value = 0.04350263745496
$ws.Range("Q3").Value = 2.016295271016667
$ws.Range("R3").Value = 18.14665743915
$ws.Range("S3").Value = 0.03565945189085799
$ws.Range("T3").Value = 0.03565945189085799

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha7"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.905615
$ws.Range("H4").Value = 35.716845
$ws.Range("I4").Value = 0.8197078149061106
$ws.Range("J4").Value = 0.8197078149061106
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.570461333333333
$ws.Range("N4").Value = 10.711384
$ws.Range("O4").Value = 0.9171442021628107
$ws.Range("P4").Value = 0.9171442021628107
$ws.Range("Q4").Value = 42.50853800705332
$ws.Range("R4").Value = 382.5768420634799
$ws.Range("S4").Value = 0.7517902699086857
$ws.Range("T4").Value = 0.7517902699086857

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha7"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.63134
$ws.Range("H5").Value = 4.894019999999999
$ws.Range("I5").Value = 0.1123186115768849
$ws.Range("J5").Value = 0.1123186115768849
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1532026666666667
$ws.Range("N5").Value = 0.459608
$ws.Range("O5").Value = 0.03935316038222934
$ws.Range("P5").Value = 0.03935316038222934
$ws.Range("Q5").Value = 0.24992563824
$ws.Range("R5").Value = 2.24933074416
$ws.Range("S5").Value = 0.004420092335294471
$ws.Range("T5").Value = 0.004420092335294471

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.63134
$ws.Range("H6").Value = 4.894019999999999
$ws.Range("I6").Value = 0.1123186115768849
$ws.Range("J6").Value = 0.1123186115768849
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1693566666666667
$ws.Range("N6").Value = 0.50807
$ws.Range("O6").Value = 0.04350263745496
$ws.Range("P6").Value = 0.04350263745496
$ws.Range("Q6").Value = 0.2762783046
$ws.Range("R6").Value = 2.4865047414
$ws.Range("S6").Value = 0.004886155838873697
$ws.Range("T6").Value = 0.004886155838873697

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha7"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.63134
$ws.Range("H7").Value = 4.894019999999999
$ws.Range("I7").Value = 0.1123186115768849
$ws.Range("J7").Value = 0.1123186115768849
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.570461333333333
$ws.Range("N7").Value = 10.711384
$ws.Range("O7").Value = 0.9171442021628107
$ws.Range("P7").Value = 0.9171442021628107
$ws.Range("Q7").Value = 5.824636391519999
$ws.Range("R7").Value = 52.42172752367999
$ws.Range("S7").Value = 0.1030123634027167
$ws.Range("T7").Value = 0.1030123634027167

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha7"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9872630000000001
$ws.Range("H8").Value = 2.961789
$ws.Range("I8").Value = 0.0679735735170045
$ws.Range("J8").Value = 0.0679735735170045
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1532026666666667
$ws.Range("N8").Value = 0.459608
$ws.Range("O8").Value = 0.03935316038222934
$ws.Range("P8").Value = 0.03935316038222934
$ws.Range("Q8").Value = 0.1512513243013334
$ws.Range("R8").Value = 1.361261918712
$ws.Range("S8").Value = 0.002674974940367935
$ws.Range("T8").Value = 0.002674974940367935

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha7"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9872630000000001
$ws.Range("H9").Value = 2.961789
$ws.Range("I9").Value = 0.0679735735170045
$ws.Range("J9").Value = 0.0679735735170045
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1693566666666667
$ws.Range("N9").Value = 0.50807
$ws.Range("O9").Value = 0.04350263745496
$ws.Range("P9").Value = 0.04350263745496
$ws.Range("Q9").Value = 0.1671995708033334
$ws.Range("R9").Value = 1.50479613723
$ws.Range("S9").Value = 0.002957029725228317
$ws.Range("T9").Value = 0.002957029725228317

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha7"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9872630000000001
$ws.Range("H10").Value = 2.961789
$ws.Range("I10").Value = 0.0679735735170045
$ws.Range("J10").Value = 0.0679735735170045
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.570461333333333
$ws.Range("N10").Value = 10.711384
$ws.Range("O10").Value = 0.9171442021628107
$ws.Range("P10").Value = 0.9171442021628107
$ws.Range("Q10").Value = 3.524984367330667
$ws.Range("R10").Value = 31.724859305976
$ws.Range("S10").Value = 0.06234156885140825
$ws.Range("T10").Value = 0.06234156885140825
